$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prakiraan Cuaca")

# Row 7 - fill in H/I/J (manual data)
$ws.Range("H7").Value = "s 10:57"
$ws.Range("I7").Value = 29
$ws.Range("J7").Value = "cerah"

# Row 8
$ws.Range("B8").Value = 31
$ws.Range("C8").Value = 64
$ws.Range("D8").Value = "Cerah"
$ws.Range("E8").Value = "E (115°)"
$ws.Range("F8").Value = 12.2
$ws.Range("G8").Value = 0

# Row 11
$ws.Range("D11").Value = "Berawan"
$ws.Range("E11").Value = "E (114°)"
$ws.Range("F11").Value = 12.9
$ws.Range("G11").Value = 0

# Row 14
$ws.Range("C14").Value = 80
$ws.Range("E14").Value = "E (102°)"
$ws.Range("F14").Value = 9.9
$ws.Range("G14").Value = 0

# Row 17
$ws.Range("C17").Value = 73
$ws.Range("D17").Value = "Berawan"
$ws.Range("E17").Value = "E (100°)"
$ws.Range("F17").Value = 17.2
$ws.Range("G17").Value = 0

# Row 20
$ws.Range("C20").Value = 66
$ws.Range("D20").Value = "Hujan Ringan"
$ws.Range("E20").Value = "SE (170°)"
$ws.Range("F20").Value = 18.5
$ws.Range("G20").Value = 1.3

# Row 23
$ws.Range("B23").Value = 30
$ws.Range("C23").Value = 72
$ws.Range("D23").Value = "Hujan Ringan"
$ws.Range("E23").Value = "E (120°)"
$ws.Range("F23").Value = 9.4
$ws.Range("G23").Value = 2

# Row 26
$ws.Range("E26").Value = "E (120°)"
$ws.Range("F26").Value = 9.2
$ws.Range("G26").Value = 0

# Row 29
$ws.Range("B29").Value = 28
$ws.Range("C29").Value = 72
$ws.Range("D29").Value = "Cerah Berawan"
$ws.Range("E29").Value = "N (40°)"
$ws.Range("F29").Value = 10.5
$ws.Range("G29").Value = 0.7

# Row 32
$ws.Range("B32").Value = 31
$ws.Range("C32").Value = 66
$ws.Range("D32").Value = "Berawan"
$ws.Range("E32").Value = "E (112°)"
$ws.Range("F32").Value = 16.7
$ws.Range("G32").Value = 0.7

# Row 35
$ws.Range("B35").Value = 30
$ws.Range("C35").Value = 71
$ws.Range("D35").Value = "Berawan"
$ws.Range("E35").Value = "E (112°)"
$ws.Range("F35").Value = 16.7
$ws.Range("G35").Value = 0

# Row 38
$ws.Range("B38").Value = 27
$ws.Range("C38").Value = 82
$ws.Range("D38").Value = "Berawan"
$ws.Range("E38").Value = "E (103°)"
$ws.Range("F38").Value = 13
$ws.Range("G38").Value = 0
